# Add three new market test-data sheets (Netherlands, Austria, Denmark) by
# copying the last sheet (Greece), which carries the right layout/styles,
# renaming it, and filling in the market-specific cells.

$wb = $excel.ActiveWorkbook

# --- Netherlands --------------------------------------------------------
$template = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $template)
$nl = $wb.Worksheets.Item($wb.Worksheets.Count)
$nl.Name = "Netherlands"
$nl.Range("B4").Value = "NGC-3144/T2176/T2177"
$nl.Range("B2").Value = "Netherlands Market"
$nl.Range("B4").Select()

# --- Austria -------------------------------------------------------------
$template = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $template)
$at = $wb.Worksheets.Item($wb.Worksheets.Count)
$at.Name = "Austria"
$at.Range("B4").Value = "NGC-3817/T2272"
$at.Range("B2").Value = "Austria Market"
$at.Range("B4").Select()

# --- Denmark ---------------------------------------------------------------
$template = $wb.Worksheets.Item($wb.Worksheets.Count)
$template.Copy($null, $template)
$dk = $wb.Worksheets.Item($wb.Worksheets.Count)
$dk.Name = "Denmark"
$dk.Range("B4").Value = "NGC-2913/T2749"
$dk.Range("B2").Value = "Denmark Market"
$dk.Range("B4").Select()

# Final tweak on Austria discovered after all three sheets were drafted.
$at.Range("A9").Value = "Fire Brigade Panel"

# Leave Austria as the active / selected sheet, scrolled so B14 is the
# last-known selection (matches the authored workbook state).
$at.Select()
$at.Range("B14").Select()

# Scroll the sheet tab strip so Swiss is the first visible tab.
$excel.ActiveWindow.ScrollWorkbookTabs(0, 3)
